# Auto commit at 2026-01-26  7:44:26.14
#
# Updates the raw metric figures on the "Metrics" sheet. The "today" sheet
# pulls every one of these through live formulas (=Metrics!Bn, =Bn, =E+...),
# so they recalculate automatically once the source cells change - no need
# to touch the formulas or their cached values by hand.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 473599
$wsMetrics.Range("B3").Value  = 362115.45000000007
$wsMetrics.Range("B4").Value  = 128146.94
$wsMetrics.Range("B5").Value  = 19281
$wsMetrics.Range("B6").Value  = 473599
$wsMetrics.Range("B7").Value  = 362115.45000000007
$wsMetrics.Range("B8").Value  = 128146.94
$wsMetrics.Range("B9").Value  = 19281
$wsMetrics.Range("B10").Value = 34574850.719999999
$wsMetrics.Range("B11").Value = 32408108.240000002
$wsMetrics.Range("B12").Value = 12073960.799999999
$wsMetrics.Range("B13").Value = 1337188

# Restore the recorded cursor position on "Metrics" (selecting a range also
# switches the active sheet, so do this before re-selecting on "today").
[void]$wsMetrics.Range("F18").Select()

# "today" is the sheet that should stay active/tabSelected afterwards, so
# select its recorded cursor position last.
$wsToday = $wb.Worksheets.Item("today")
[void]$wsToday.Range("I10").Select()
